# testRap.docx fix-up:
#  1) Drop the stray <w:proofErr gramStart/gramEnd> pair wrapping the
#     centered "à" paragraph (harmless spell/grammar-check artifacts,
#     no visible text change).
#  2) Drop the empty "charge"-styled paragraph that was left dangling
#     right after the "80 kWh x 0.15€ = 12€" line.

$d = $word.ActiveDocument

function Trim-ParaText($range) {
    return $range.Text.TrimEnd([char]13, [char]7)
}

# --- 1) Remove the proofErr wrapper around the "à" paragraph ------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ((Trim-ParaText $p.Range) -eq "à") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6FD33F53" w14:textId="075154E6" w:rsidR="007931C9" w:rsidRPr="003865A4" w:rsidRDefault="007931C9" w:rsidP="007931C9"><w:pPr><w:jc w:val="center"/><w:rPr><w:szCs w:val="20"/></w:rPr></w:pPr><w:r w:rsidRPr="003865A4"><w:rPr><w:szCs w:val="20"/></w:rPr><w:t>à</w:t></w:r></w:p>'
    [void]$target.Range.InsertXML($xml)
}

# --- 2) Remove the dangling empty "charge" paragraph after the kWh line -
# (select from the end of the previous paragraph through the end of this
#  one so the whole paragraph - including its own mark - collapses away;
#  deleting only Paragraphs(i).Range here leaves the pilcrow behind in
#  the saved package even though the live Paragraphs collection looks
#  right, so we fold the mark into the preceding paragraph instead.)
for ($i = 2; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $prev = $d.Paragraphs($i - 1)
    if ((Trim-ParaText $p.Range) -eq "" -and
        (Trim-ParaText $prev.Range) -eq "80 kWh x 0.15€ = 12€") {
        $delRange = $d.Range($p.Range.Start - 1, $p.Range.End)
        $delRange.Delete()
        break
    }
}
